# OrangeHRM_test_cases.xlsx - add new Recruitment/Candidates test case (row 11)
# and restyle the following blank rows (14-20) to match the bordered block
# used by rows 4-13.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1) New test-case row content (row 11)
# ------------------------------------------------------------------
$ws.Range("A11").Value = "TS_OHRM_002"
$ws.Range("B11").Value = "Verify Recruitment page Candidates tab functionality "
$ws.Range("C11").Value = "TC_OHRM_RPCT_001"
$ws.Range("D11").Value = "Add candidate with correct required data only"
$ws.Range("F11").Value = "FirstName - ""Stan""`nLastName - ""Smith"" `nEmail - ""stan_smith@gmail.com"""
$ws.Range("E11").Value = "1. Successfully login to the app`n2. Click on 'Recruitment' menu on the side pannel`n3. Verify location is 'Candidates' tab`n4. Verify there is no record with Test Data details in the 'Records Found' section on the page`n5. Click 'Add' button`n6. On the appeared form input Test Data (FirstName, LastName, Email) and click 'Save' button.`n7. Verify record with Test Data details appeared on the 'Records Found' section on the page."
$ws.Range("G11").Value = "Success' green toast message should appear and Candidate profile should be on the screen after saving.`nThe record with Test Data details should appear on the 'Records Found' section on the page."

# Row 11 grows tall to fit the wrapped multi-line text.
$ws.Rows.Item(11).RowHeight = 150

# ------------------------------------------------------------------
# 2) G11 needs the "quote prefix" flavour of the wrap/border style
#    (Excel sets this automatically whenever a cell's typed text starts
#    with an apostrophe). Build it on a scratch cell, copy the format
#    across, then clean the scratch cell up again.
# ------------------------------------------------------------------
$ws.Range("E4").Copy()
$ws.Range("Z1").PasteSpecial(-4122)
$ws.Range("Z1").Value = "'x"
$ws.Range("Z1").Copy()
$ws.Range("G11").PasteSpecial(-4122)
$ws.Range("Z1").Clear()

# ------------------------------------------------------------------
# 3) Rows 14-20 switch from the plain style to the bordered/top-aligned
#    style already used by rows 4-13 (A/H columns vs B:G columns).
# ------------------------------------------------------------------
$ws.Range("A13").Copy()
$ws.Range("A14:A20").PasteSpecial(-4122)

$ws.Range("B13:G13").Copy()
$ws.Range("B14:G20").PasteSpecial(-4122)

$ws.Range("H13").Copy()
$ws.Range("H14:H20").PasteSpecial(-4122)

# ------------------------------------------------------------------
# 4) View state: scroll so row 9 is pinned below the frozen header and
#    column B is the first visible column, with G12 as the active cell.
# ------------------------------------------------------------------
$ws.Range("G12").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 9
$win.ScrollColumn = 2
